$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (FiT) and C (Price) then D/E/F recalculated values
$ws.Range("B2:B8").Value = 0.03
$ws.Range("C2:C8").Value = 0.28999999999999998

$dValues = @(0.083427973557962498, 0.08342797355796247, 0.083427973557962457, 0.083427973557949911, 0.083427973557954088, 0.083427973557961096, 0.083427973557962262)
$eValues = @(0.62994710897297823, 0.62994710897297856, 0.62994710897297856, 0.62994710897297701, 0.62994710897297823, 0.62994710897297579, 0.62994710897297779)
$fValues = @(59533811.581426457, 59533811.581426471, 59533811.581426457, 59533811.581426948, 59533811.581426769, 59533811.581426539, 59533811.581426486)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
    $ws.Cells.Item($row, 6).Value = $fValues[$i]
}

$ws.Range("D2:E8").Style = "Percent"
$ws.Range("F2:F8").Style = "Currency"

# Best-fit width on column F (closest attainable value to the recorded
# 14.6328125 "best fit" width given this engine's column-width quantisation)
$ws.Columns.Item(6).ColumnWidth = 13.9

# Matches the saved selection (active cell) recorded in the workbook
$ws.Range("F3").Select() | Out-Null
